$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2181208053691275
$ws.Range("C2").Value = 0.5167785234899329
$ws.Range("J2").Value = 0.01006711409395973
$ws.Range("P2").Value = 0.1644295302013423
$ws.Range("S2").Value = 0.09060402684563758
$ws.Range("B3").Value = 0.006060606060606061
$ws.Range("C3").Value = 0.06666666666666667
$ws.Range("J3").Value = 0.04242424242424243
$ws.Range("P3").Value = 0.6848484848484848
$ws.Range("S3").Value = 0.2
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.4
$ws.Range("B6").Value = 0.08292682926829269
$ws.Range("F6").Value = 0.04390243902439024
$ws.Range("J6").Value = 0.2634146341463415
$ws.Range("O6").Value = 0.03902439024390244
$ws.Range("Q6").Value = 0.1609756097560976
$ws.Range("R6").Value = 0.03902439024390244
$ws.Range("S6").Value = 0.3707317073170732
$ws.Range("B7").Value = 0.0949367088607595
$ws.Range("D7").Value = 0.03164556962025317
$ws.Range("E7").Value = 0.006329113924050633
$ws.Range("F7").Value = 0.04430379746835443
$ws.Range("J7").Value = 0.1265822784810127
$ws.Range("O7").Value = 0.0189873417721519
$ws.Range("Q7").Value = 0.1708860759493671
$ws.Range("R7").Value = 0.0379746835443038
$ws.Range("S7").Value = 0.4683544303797468
$ws.Range("B8").Value = 0.1005025125628141
$ws.Range("D8").Value = 0.03517587939698492
$ws.Range("E8").Value = 0.002512562814070352
$ws.Range("F8").Value = 0.05025125628140704
$ws.Range("J8").Value = 0.08542713567839195
$ws.Range("O8").Value = 0.01507537688442211
$ws.Range("Q8").Value = 0.1884422110552764
$ws.Range("R8").Value = 0.1005025125628141
$ws.Range("S8").Value = 0.4221105527638191
$ws.Range("B9").Value = 0.102803738317757
$ws.Range("D9").Value = 0.01401869158878505
$ws.Range("F9").Value = 0.04672897196261682
$ws.Range("J9").Value = 0.08878504672897196
$ws.Range("O9").Value = 0.01401869158878505
$ws.Range("Q9").Value = 0.2383177570093458
$ws.Range("R9").Value = 0.07009345794392523
$ws.Range("S9").Value = 0.4252336448598131
$ws.Range("B10").Value = 0.1068273092369478
$ws.Range("D10").Value = 0.02248995983935743
$ws.Range("E10").Value = 0.002409638554216868
$ws.Range("F10").Value = 0.07228915662650602
$ws.Range("J10").Value = 0.1092369477911647
$ws.Range("O10").Value = 0.01686746987951807
$ws.Range("Q10").Value = 0.1975903614457831
$ws.Range("R10").Value = 0.08995983935742972
$ws.Range("S10").Value = 0.3823293172690763
$ws.Range("G11").Value = 0.1433447098976109
$ws.Range("J11").Value = 0.1262798634812287
$ws.Range("K11").Value = 0.2184300341296928
$ws.Range("L11").Value = 0.4948805460750853
$ws.Range("S11").Value = 0.01706484641638225
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.2585034013605442
$ws.Range("L12").Value = 0.01360544217687075
$ws.Range("S12").Value = 0.06122448979591837
$ws.Range("G13").Value = 0.5365853658536586
$ws.Range("J13").Value = 0.4146341463414634
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.01382488479262673
$ws.Range("H15").Value = 0.1382488479262673
$ws.Range("I15").Value = 0.08294930875576037
$ws.Range("J15").Value = 0.3732718894009217
$ws.Range("K15").Value = 0.07834101382488479
$ws.Range("M15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.05990783410138249
$ws.Range("S15").Value = 0.2488479262672811
$ws.Range("F16").Value = 0.02150537634408602
$ws.Range("H16").Value = 0.1182795698924731
$ws.Range("I16").Value = 0.09677419354838709
$ws.Range("J16").Value = 0.4516129032258064
$ws.Range("K16").Value = 0.08064516129032258
$ws.Range("M16").Value = 0.02150537634408602
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.1451612903225807
$ws.Range("F17").Value = 0.02097902097902098
$ws.Range("H17").Value = 0.1748251748251748
$ws.Range("I17").Value = 0.1142191142191142
$ws.Range("J17").Value = 0.4102564102564102
$ws.Range("K17").Value = 0.09090909090909091
$ws.Range("M17").Value = 0.01398601398601399
$ws.Range("O17").Value = 0.04895104895104895
$ws.Range("S17").Value = 0.1258741258741259
$ws.Range("F18").Value = 0.01666666666666667
$ws.Range("H18").Value = 0.1833333333333333
$ws.Range("I18").Value = 0.08888888888888889
$ws.Range("J18").Value = 0.3555555555555556
$ws.Range("K18").Value = 0.08333333333333333
$ws.Range("M18").Value = 0.01666666666666667
$ws.Range("O18").Value = 0.07777777777777778
$ws.Range("S18").Value = 0.1777777777777778
$ws.Range("F19").Value = 0.01381427475057559
$ws.Range("H19").Value = 0.1849577897160399
$ws.Range("I19").Value = 0.08902532617037605
$ws.Range("J19").Value = 0.3737528779739064
$ws.Range("K19").Value = 0.1105141980046048
$ws.Range("M19").Value = 0.02072141212586339
$ws.Range("N19").Value = 0.0007674597083653108
$ws.Range("O19").Value = 0.07444359171143515
$ws.Range("S19").Value = 0.1320030698388334
